$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("runs")
$ws1.Range("M2").Value = 45.50528526306152
$ws1.Range("N2").Value = 45.50528526306152
$ws1.Range("L3").Value = 4.999637603759766
$ws1.Range("M3").Value = 1.981973648071289
$ws1.Range("N3").Value = 6.981611251831055
$ws1.Range("M4").Value = 37.02330589294434
$ws1.Range("N4").Value = 37.02330589294434
$ws1.Range("L5").Value = 16.00003242492676
$ws1.Range("M5").Value = 16.00313186645508
$ws1.Range("N5").Value = 32.00316429138184
$ws1.Range("M6").Value = 34.69443321228027
$ws1.Range("N6").Value = 34.69443321228027
$ws1.Range("L7").Value = 14.21213150024414
$ws1.Range("M7").Value = 18.9812183380127
$ws1.Range("N7").Value = 33.19334983825684
$ws1.Range("M8").Value = 33.36477279663086
$ws1.Range("N8").Value = 33.36477279663086
$ws1.Range("L9").Value = 19.18649673461914
$ws1.Range("M9").Value = 23.53835105895996
$ws1.Range("N9").Value = 42.7248477935791
$ws1.Range("M10").Value = 194.9613094329834
$ws1.Range("N10").Value = 194.9613094329834
$ws1.Range("L11").Value = 20.22838592529297
$ws1.Range("M11").Value = 64.35441970825195
$ws1.Range("N11").Value = 84.58280563354492
$ws1.Range("M12").Value = 38.47455978393555
$ws1.Range("N12").Value = 38.47455978393555
$ws1.Range("L13").Value = 4.00090217590332
$ws1.Range("M13").Value = 3.000974655151367
$ws1.Range("N13").Value = 7.001876831054688
$ws1.Range("M14").Value = 14.00160789489746
$ws1.Range("N14").Value = 14.00160789489746
$ws1.Range("L15").Value = 5.048274993896484
$ws1.Range("M15").Value = 3.061056137084961
$ws1.Range("N15").Value = 8.109331130981445
$ws1.Range("M16").Value = 12.99929618835449
$ws1.Range("N16").Value = 12.99929618835449
$ws1.Range("L17").Value = 6.430149078369141
$ws1.Range("M17").Value = 4.999399185180664
$ws1.Range("N17").Value = 11.4295482635498
$ws1.Range("M18").Value = 10.32590866088867
$ws1.Range("N18").Value = 10.32590866088867
$ws1.Range("L19").Value = 18.97525787353516
$ws1.Range("M19").Value = 26.0016918182373
$ws1.Range("N19").Value = 44.97694969177246
$ws1.Range("M20").Value = 15.99645614624023
$ws1.Range("N20").Value = 15.99645614624023
$ws1.Range("L21").Value = 19.00124549865723
$ws1.Range("M21").Value = 17.85397529602051
$ws1.Range("N21").Value = 36.85522079467773
$ws1.Range("M22").Value = 5.548238754272461
$ws1.Range("N22").Value = 5.548238754272461
$ws1.Range("L23").Value = 4.98652458190918
$ws1.Range("M23").Value = 0.9953975677490234
$ws1.Range("N23").Value = 5.981922149658203
$ws1.Range("M24").Value = 4.981517791748047
$ws1.Range("N24").Value = 4.981517791748047
$ws1.Range("L25").Value = 3.97944450378418
$ws1.Range("M25").Value = 1.009702682495117
$ws1.Range("N25").Value = 4.989147186279297
$ws1.Range("M26").Value = 5.998134613037109
$ws1.Range("N26").Value = 5.998134613037109
$ws1.Range("L27").Value = 6.001710891723633
$ws1.Range("M27").Value = 2.002716064453125
$ws1.Range("N27").Value = 8.004426956176758
$ws1.Range("M28").Value = 7.02357292175293
$ws1.Range("N28").Value = 7.02357292175293
$ws1.Range("L29").Value = 17.50564575195312
$ws1.Range("M29").Value = 1.001358032226562
$ws1.Range("N29").Value = 18.50700378417969
$ws1.Range("M30").Value = 6.006479263305664
$ws1.Range("N30").Value = 6.006479263305664
$ws1.Range("L31").Value = 15.85698127746582
$ws1.Range("M31").Value = 2.000331878662109
$ws1.Range("N31").Value = 17.85731315612793
$ws1.Range("M32").Value = 32.52673149108887
$ws1.Range("N32").Value = 32.52673149108887
$ws1.Range("L33").Value = 5.06901741027832
$ws1.Range("M33").Value = 0.9996891021728516
$ws1.Range("N33").Value = 6.068706512451172
$ws1.Range("M34").Value = 24.23977851867676
$ws1.Range("N34").Value = 24.23977851867676
$ws1.Range("L35").Value = 5.043268203735352
$ws1.Range("M35").Value = 1.028060913085938
$ws1.Range("N35").Value = 6.071329116821289
$ws1.Range("M36").Value = 44.92068290710449
$ws1.Range("N36").Value = 44.92068290710449
$ws1.Range("L37").Value = 5.000591278076172
$ws1.Range("M37").Value = 1.997709274291992
$ws1.Range("N37").Value = 6.998300552368164
$ws1.Range("M38").Value = 20.45631408691406
$ws1.Range("N38").Value = 20.45631408691406
$ws1.Range("L39").Value = 5.982160568237305
$ws1.Range("M39").Value = 1.062393188476562
$ws1.Range("N39").Value = 7.044553756713867
$ws1.Range("M40").Value = 18.93782615661621
$ws1.Range("N40").Value = 18.93782615661621
$ws1.Range("L41").Value = 17.73238182067871
$ws1.Range("M41").Value = 2.00963020324707
$ws1.Range("N41").Value = 19.74201202392578
$ws1.Range("M42").Value = 52.1690845489502
$ws1.Range("N42").Value = 52.1690845489502
$ws1.Range("L43").Value = 6.330013275146484
$ws1.Range("M43").Value = 0.9987354278564453
$ws1.Range("N43").Value = 7.32874870300293
$ws1.Range("M44").Value = 43.59340667724609
$ws1.Range("N44").Value = 43.59340667724609
$ws1.Range("L45").Value = 5.985498428344727
$ws1.Range("M45").Value = 1.000165939331055
$ws1.Range("N45").Value = 6.985664367675781
$ws1.Range("M46").Value = 55.08303642272949
$ws1.Range("N46").Value = 55.08303642272949
$ws1.Range("L47").Value = 4.999876022338867
$ws1.Range("M47").Value = 1.003503799438477
$ws1.Range("N47").Value = 6.003379821777344
$ws1.Range("M48").Value = 41.64719581604004
$ws1.Range("N48").Value = 41.64719581604004
$ws1.Range("L49").Value = 4.147529602050781
$ws1.Range("M49").Value = 1.507997512817383
$ws1.Range("N49").Value = 5.655527114868164
$ws1.Range("M50").Value = 54.79264259338379
$ws1.Range("N50").Value = 54.79264259338379
$ws1.Range("L51").Value = 20.13683319091797
$ws1.Range("M51").Value = 0.965118408203125
$ws1.Range("N51").Value = 21.10195159912109
$ws1.Range("M52").Value = 46.54788970947266
$ws1.Range("N52").Value = 46.54788970947266
$ws1.Range("L53").Value = 5.097389221191406
$ws1.Range("M53").Value = 1.999616622924805
$ws1.Range("N53").Value = 7.097005844116211
$ws1.Range("M54").Value = 12.80665397644043
$ws1.Range("N54").Value = 12.80665397644043
$ws1.Range("L55").Value = 7.826089859008789
$ws1.Range("M55").Value = 4.071474075317383
$ws1.Range("N55").Value = 11.89756393432617
$ws1.Range("M56").Value = 17.9448127746582
$ws1.Range("N56").Value = 17.9448127746582
$ws1.Range("L57").Value = 4.979133605957031
$ws1.Range("M57").Value = 3.999233245849609
$ws1.Range("N57").Value = 8.978366851806641
$ws1.Range("M58").Value = 12.21728324890137
$ws1.Range("N58").Value = 12.21728324890137
$ws1.Range("L59").Value = 9.620428085327148
$ws1.Range("M59").Value = 9.126901626586914
$ws1.Range("N59").Value = 18.74732971191406
$ws1.Range("M60").Value = 12.8328800201416
$ws1.Range("N60").Value = 12.8328800201416
$ws1.Range("L61").Value = 18.77355575561523
$ws1.Range("M61").Value = 18.14031600952148
$ws1.Range("N61").Value = 36.91387176513672
$ws1.Range("M62").Value = 23.15688133239746
$ws1.Range("N62").Value = 23.15688133239746
$ws1.Range("L63").Value = 3.980875015258789
$ws1.Range("N63").Value = 3.980875015258789
$ws1.Range("M64").Value = 22.99213409423828
$ws1.Range("N64").Value = 22.99213409423828
$ws1.Range("L65").Value = 7.234811782836914
$ws1.Range("N65").Value = 7.234811782836914
$ws1.Range("M66").Value = 54.16202545166016
$ws1.Range("N66").Value = 54.16202545166016
$ws1.Range("L67").Value = 10.50329208374023
$ws1.Range("N67").Value = 10.50329208374023
$ws1.Range("M68").Value = 34.48224067687988
$ws1.Range("N68").Value = 34.48224067687988
$ws1.Range("L69").Value = 9.414196014404297
$ws1.Range("N69").Value = 9.414196014404297
$ws1.Range("M70").Value = 198.979377746582
$ws1.Range("N70").Value = 198.979377746582
$ws1.Range("L71").Value = 17.48251914978027
$ws1.Range("M71").Value = 58.32958221435547
$ws1.Range("N71").Value = 75.81210136413574
$ws1.Range("M72").Value = 35.49695014953613
$ws1.Range("N72").Value = 35.49695014953613
$ws1.Range("L73").Value = 5.006551742553711
$ws1.Range("M73").Value = 2.999782562255859
$ws1.Range("N73").Value = 8.00633430480957
$ws1.Range("M74").Value = 43.0457592010498
$ws1.Range("N74").Value = 43.0457592010498
$ws1.Range("L75").Value = 4.985332489013672
$ws1.Range("M75").Value = 4.049539566040039
$ws1.Range("N75").Value = 9.034872055053711
$ws1.Range("M76").Value = 55.21059036254883
$ws1.Range("N76").Value = 55.21059036254883
$ws1.Range("L77").Value = 4.991292953491211
$ws1.Range("M77").Value = 5.012035369873047
$ws1.Range("N77").Value = 10.00332832336426
$ws1.Range("M78").Value = 10.00785827636719
$ws1.Range("N78").Value = 10.00785827636719
$ws1.Range("L79").Value = 12.98832893371582
$ws1.Range("N79").Value = 12.98832893371582
$ws1.Range("M80").Value = 5.960226058959961
$ws1.Range("N80").Value = 5.960226058959961
$ws1.Range("L81").Value = 15.98763465881348
$ws1.Range("M81").Value = 1.000404357910156
$ws1.Range("N81").Value = 16.98803901672363
$ws1.Range("M82").Value = 39.22796249389648
$ws1.Range("N82").Value = 39.22796249389648
$ws1.Range("L83").Value = 4.997491836547852
$ws1.Range("M83").Value = 2.999305725097656
$ws1.Range("N83").Value = 7.996797561645508
$ws1.Range("M84").Value = 35.37106513977051
$ws1.Range("N84").Value = 35.37106513977051
$ws1.Range("L85").Value = 5.545377731323242
$ws1.Range("M85").Value = 2.999067306518555
$ws1.Range("N85").Value = 8.544445037841797
$ws1.Range("M86").Value = 21.99888229370117
$ws1.Range("N86").Value = 21.99888229370117
$ws1.Range("L87").Value = 7.001638412475586
$ws1.Range("N87").Value = 7.001638412475586
$ws1.Range("M88").Value = 34.03639793395996
$ws1.Range("N88").Value = 34.03639793395996
$ws1.Range("L89").Value = 20.01690864562988
$ws1.Range("N89").Value = 20.01690864562988
$ws1.Range("L90").Value = 0
$ws1.Range("M90").Value = 28.9008617401123
$ws1.Range("N90").Value = 28.9008617401123
$ws1.Range("L91").Value = 22.11570739746094
$ws1.Range("M91").Value = 3.004074096679688
$ws1.Range("N91").Value = 25.11978149414062
$ws1.Range("M92").Value = 31.16607666015625
$ws1.Range("N92").Value = 31.16607666015625
$ws1.Range("L93").Value = 3.986120223999023
$ws1.Range("M93").Value = 3.000020980834961
$ws1.Range("N93").Value = 6.986141204833984
$ws1.Range("M94").Value = 29.54363822937012
$ws1.Range("N94").Value = 29.54363822937012
$ws1.Range("L95").Value = 6.667375564575195
$ws1.Range("M95").Value = 5.9967041015625
$ws1.Range("N95").Value = 12.6640796661377
$ws1.Range("M96").Value = 20.58196067810059
$ws1.Range("N96").Value = 20.58196067810059
$ws1.Range("L97").Value = 9.058713912963867
$ws1.Range("M97").Value = 9.058475494384766
$ws1.Range("N97").Value = 18.11718940734863
$ws1.Range("M98").Value = 19.93465423583984
$ws1.Range("N98").Value = 19.93465423583984
$ws1.Range("L99").Value = 17.22240447998047
$ws1.Range("M99").Value = 27.70376205444336
$ws1.Range("N99").Value = 44.92616653442383
$ws1.Range("M100").Value = 187.4191761016846
$ws1.Range("N100").Value = 187.4191761016846
$ws1.Range("L101").Value = 24.48391914367676
$ws1.Range("M101").Value = 67.46840476989746
$ws1.Range("N101").Value = 91.95232391357422
$ws2 = $wb.Worksheets.Item("runs_flat")
$ws2.Range("L2").Value = 4.999637603759766
$ws2.Range("M2").Value = 1.981973648071289
$ws2.Range("N2").Value = 6.981611251831055
$ws2.Range("L3").Value = 4.999637603759766
$ws2.Range("M3").Value = 1.981973648071289
$ws2.Range("N3").Value = 6.981611251831055
$ws2.Range("L4").Value = 16.00003242492676
$ws2.Range("M4").Value = 16.00313186645508
$ws2.Range("N4").Value = 32.00316429138184
$ws2.Range("L5").Value = 16.00003242492676
$ws2.Range("M5").Value = 16.00313186645508
$ws2.Range("N5").Value = 32.00316429138184
$ws2.Range("L6").Value = 14.21213150024414
$ws2.Range("M6").Value = 18.9812183380127
$ws2.Range("N6").Value = 33.19334983825684
$ws2.Range("L7").Value = 14.21213150024414
$ws2.Range("M7").Value = 18.9812183380127
$ws2.Range("N7").Value = 33.19334983825684
$ws2.Range("L8").Value = 19.18649673461914
$ws2.Range("M8").Value = 23.53835105895996
$ws2.Range("N8").Value = 42.7248477935791
$ws2.Range("L9").Value = 19.18649673461914
$ws2.Range("M9").Value = 23.53835105895996
$ws2.Range("N9").Value = 42.7248477935791
$ws2.Range("M10").Value = 194.9613094329834
$ws2.Range("N10").Value = 194.9613094329834
$ws2.Range("M11").Value = 194.9613094329834
$ws2.Range("N11").Value = 194.9613094329834
$ws2.Range("M12").Value = 194.9613094329834
$ws2.Range("N12").Value = 194.9613094329834
$ws2.Range("M13").Value = 194.9613094329834
$ws2.Range("N13").Value = 194.9613094329834
$ws2.Range("M14").Value = 194.9613094329834
$ws2.Range("N14").Value = 194.9613094329834
$ws2.Range("M15").Value = 194.9613094329834
$ws2.Range("N15").Value = 194.9613094329834
$ws2.Range("L16").Value = 20.22838592529297
$ws2.Range("M16").Value = 64.35441970825195
$ws2.Range("N16").Value = 84.58280563354492
$ws2.Range("L17").Value = 20.22838592529297
$ws2.Range("M17").Value = 64.35441970825195
$ws2.Range("N17").Value = 84.58280563354492
$ws2.Range("L18").Value = 20.22838592529297
$ws2.Range("M18").Value = 64.35441970825195
$ws2.Range("N18").Value = 84.58280563354492
$ws2.Range("L19").Value = 20.22838592529297
$ws2.Range("M19").Value = 64.35441970825195
$ws2.Range("N19").Value = 84.58280563354492
$ws2.Range("L20").Value = 4.00090217590332
$ws2.Range("M20").Value = 3.000974655151367
$ws2.Range("N20").Value = 7.001876831054688
$ws2.Range("L21").Value = 4.00090217590332
$ws2.Range("M21").Value = 3.000974655151367
$ws2.Range("N21").Value = 7.001876831054688
$ws2.Range("L22").Value = 5.048274993896484
$ws2.Range("M22").Value = 3.061056137084961
$ws2.Range("N22").Value = 8.109331130981445
$ws2.Range("L23").Value = 5.048274993896484
$ws2.Range("M23").Value = 3.061056137084961
$ws2.Range("N23").Value = 8.109331130981445
$ws2.Range("L24").Value = 6.430149078369141
$ws2.Range("M24").Value = 4.999399185180664
$ws2.Range("N24").Value = 11.4295482635498
$ws2.Range("L25").Value = 6.430149078369141
$ws2.Range("M25").Value = 4.999399185180664
$ws2.Range("N25").Value = 11.4295482635498
$ws2.Range("L26").Value = 18.97525787353516
$ws2.Range("M26").Value = 26.0016918182373
$ws2.Range("N26").Value = 44.97694969177246
$ws2.Range("L27").Value = 18.97525787353516
$ws2.Range("M27").Value = 26.0016918182373
$ws2.Range("N27").Value = 44.97694969177246
$ws2.Range("L28").Value = 19.00124549865723
$ws2.Range("M28").Value = 17.85397529602051
$ws2.Range("N28").Value = 36.85522079467773
$ws2.Range("L29").Value = 19.00124549865723
$ws2.Range("M29").Value = 17.85397529602051
$ws2.Range("N29").Value = 36.85522079467773
$ws2.Range("M30").Value = 5.548238754272461
$ws2.Range("N30").Value = 5.548238754272461
$ws2.Range("L31").Value = 4.98652458190918
$ws2.Range("M31").Value = 0.9953975677490234
$ws2.Range("N31").Value = 5.981922149658203
$ws2.Range("M32").Value = 4.981517791748047
$ws2.Range("N32").Value = 4.981517791748047
$ws2.Range("L33").Value = 3.97944450378418
$ws2.Range("M33").Value = 1.009702682495117
$ws2.Range("N33").Value = 4.989147186279297
$ws2.Range("M34").Value = 5.998134613037109
$ws2.Range("N34").Value = 5.998134613037109
$ws2.Range("L35").Value = 6.001710891723633
$ws2.Range("M35").Value = 2.002716064453125
$ws2.Range("N35").Value = 8.004426956176758
$ws2.Range("M36").Value = 7.02357292175293
$ws2.Range("N36").Value = 7.02357292175293
$ws2.Range("L37").Value = 17.50564575195312
$ws2.Range("M37").Value = 1.001358032226562
$ws2.Range("N37").Value = 18.50700378417969
$ws2.Range("M38").Value = 6.006479263305664
$ws2.Range("N38").Value = 6.006479263305664
$ws2.Range("L39").Value = 15.85698127746582
$ws2.Range("M39").Value = 2.000331878662109
$ws2.Range("N39").Value = 17.85731315612793
$ws2.Range("M40").Value = 32.52673149108887
$ws2.Range("N40").Value = 32.52673149108887
$ws2.Range("M41").Value = 32.52673149108887
$ws2.Range("N41").Value = 32.52673149108887
$ws2.Range("L42").Value = 5.06901741027832
$ws2.Range("M42").Value = 0.9996891021728516
$ws2.Range("N42").Value = 6.068706512451172
$ws2.Range("M43").Value = 24.23977851867676
$ws2.Range("N43").Value = 24.23977851867676
$ws2.Range("L44").Value = 5.043268203735352
$ws2.Range("M44").Value = 1.028060913085938
$ws2.Range("N44").Value = 6.071329116821289
$ws2.Range("M45").Value = 44.92068290710449
$ws2.Range("N45").Value = 44.92068290710449
$ws2.Range("M46").Value = 44.92068290710449
$ws2.Range("N46").Value = 44.92068290710449
$ws2.Range("L47").Value = 5.000591278076172
$ws2.Range("M47").Value = 1.997709274291992
$ws2.Range("N47").Value = 6.998300552368164
$ws2.Range("M48").Value = 20.45631408691406
$ws2.Range("N48").Value = 20.45631408691406
$ws2.Range("L49").Value = 5.982160568237305
$ws2.Range("M49").Value = 1.062393188476562
$ws2.Range("N49").Value = 7.044553756713867
$ws2.Range("M50").Value = 18.93782615661621
$ws2.Range("N50").Value = 18.93782615661621
$ws2.Range("L51").Value = 17.73238182067871
$ws2.Range("M51").Value = 2.00963020324707
$ws2.Range("N51").Value = 19.74201202392578
$ws2.Range("M52").Value = 52.1690845489502
$ws2.Range("N52").Value = 52.1690845489502
$ws2.Range("L53").Value = 6.330013275146484
$ws2.Range("M53").Value = 0.9987354278564453
$ws2.Range("N53").Value = 7.32874870300293
$ws2.Range("M54").Value = 43.59340667724609
$ws2.Range("N54").Value = 43.59340667724609
$ws2.Range("L55").Value = 5.985498428344727
$ws2.Range("M55").Value = 1.000165939331055
$ws2.Range("N55").Value = 6.985664367675781
$ws2.Range("M56").Value = 55.08303642272949
$ws2.Range("N56").Value = 55.08303642272949
$ws2.Range("L57").Value = 4.999876022338867
$ws2.Range("M57").Value = 1.003503799438477
$ws2.Range("N57").Value = 6.003379821777344
$ws2.Range("M58").Value = 41.64719581604004
$ws2.Range("N58").Value = 41.64719581604004
$ws2.Range("L59").Value = 4.147529602050781
$ws2.Range("M59").Value = 1.507997512817383
$ws2.Range("N59").Value = 5.655527114868164
$ws2.Range("M60").Value = 54.79264259338379
$ws2.Range("N60").Value = 54.79264259338379
$ws2.Range("L61").Value = 20.13683319091797
$ws2.Range("M61").Value = 0.965118408203125
$ws2.Range("N61").Value = 21.10195159912109
$ws2.Range("L62").Value = 5.097389221191406
$ws2.Range("M62").Value = 1.999616622924805
$ws2.Range("N62").Value = 7.097005844116211
$ws2.Range("L63").Value = 5.097389221191406
$ws2.Range("M63").Value = 1.999616622924805
$ws2.Range("N63").Value = 7.097005844116211
$ws2.Range("L64").Value = 7.826089859008789
$ws2.Range("M64").Value = 4.071474075317383
$ws2.Range("N64").Value = 11.89756393432617
$ws2.Range("L65").Value = 7.826089859008789
$ws2.Range("M65").Value = 4.071474075317383
$ws2.Range("N65").Value = 11.89756393432617
$ws2.Range("L66").Value = 4.979133605957031
$ws2.Range("M66").Value = 3.999233245849609
$ws2.Range("N66").Value = 8.978366851806641
$ws2.Range("L67").Value = 4.979133605957031
$ws2.Range("M67").Value = 3.999233245849609
$ws2.Range("N67").Value = 8.978366851806641
$ws2.Range("L68").Value = 9.620428085327148
$ws2.Range("M68").Value = 9.126901626586914
$ws2.Range("N68").Value = 18.74732971191406
$ws2.Range("L69").Value = 9.620428085327148
$ws2.Range("M69").Value = 9.126901626586914
$ws2.Range("N69").Value = 18.74732971191406
$ws2.Range("L70").Value = 18.77355575561523
$ws2.Range("M70").Value = 18.14031600952148
$ws2.Range("N70").Value = 36.91387176513672
$ws2.Range("L71").Value = 18.77355575561523
$ws2.Range("M71").Value = 18.14031600952148
$ws2.Range("N71").Value = 36.91387176513672
$ws2.Range("M72").Value = 23.15688133239746
$ws2.Range("N72").Value = 23.15688133239746
$ws2.Range("L73").Value = 3.980875015258789
$ws2.Range("N73").Value = 3.980875015258789
$ws2.Range("M74").Value = 22.99213409423828
$ws2.Range("N74").Value = 22.99213409423828
$ws2.Range("L75").Value = 7.234811782836914
$ws2.Range("N75").Value = 7.234811782836914
$ws2.Range("L76").Value = 10.50329208374023
$ws2.Range("N76").Value = 10.50329208374023
$ws2.Range("M77").Value = 34.48224067687988
$ws2.Range("N77").Value = 34.48224067687988
$ws2.Range("L78").Value = 9.414196014404297
$ws2.Range("N78").Value = 9.414196014404297
$ws2.Range("M79").Value = 198.979377746582
$ws2.Range("N79").Value = 198.979377746582
$ws2.Range("M80").Value = 198.979377746582
$ws2.Range("N80").Value = 198.979377746582
$ws2.Range("M81").Value = 198.979377746582
$ws2.Range("N81").Value = 198.979377746582
$ws2.Range("M82").Value = 198.979377746582
$ws2.Range("N82").Value = 198.979377746582
$ws2.Range("M83").Value = 198.979377746582
$ws2.Range("N83").Value = 198.979377746582
$ws2.Range("M84").Value = 198.979377746582
$ws2.Range("N84").Value = 198.979377746582
$ws2.Range("L85").Value = 17.48251914978027
$ws2.Range("M85").Value = 58.32958221435547
$ws2.Range("N85").Value = 75.81210136413574
$ws2.Range("L86").Value = 17.48251914978027
$ws2.Range("M86").Value = 58.32958221435547
$ws2.Range("N86").Value = 75.81210136413574
$ws2.Range("L87").Value = 17.48251914978027
$ws2.Range("M87").Value = 58.32958221435547
$ws2.Range("N87").Value = 75.81210136413574
$ws2.Range("L88").Value = 17.48251914978027
$ws2.Range("M88").Value = 58.32958221435547
$ws2.Range("N88").Value = 75.81210136413574
$ws2.Range("M89").Value = 35.49695014953613
$ws2.Range("N89").Value = 35.49695014953613
$ws2.Range("M90").Value = 35.49695014953613
$ws2.Range("N90").Value = 35.49695014953613
$ws2.Range("L91").Value = 5.006551742553711
$ws2.Range("M91").Value = 2.999782562255859
$ws2.Range("N91").Value = 8.00633430480957
$ws2.Range("L92").Value = 5.006551742553711
$ws2.Range("M92").Value = 2.999782562255859
$ws2.Range("N92").Value = 8.00633430480957
$ws2.Range("M93").Value = 43.0457592010498
$ws2.Range("N93").Value = 43.0457592010498
$ws2.Range("M94").Value = 43.0457592010498
$ws2.Range("N94").Value = 43.0457592010498
$ws2.Range("L95").Value = 4.985332489013672
$ws2.Range("M95").Value = 4.049539566040039
$ws2.Range("N95").Value = 9.034872055053711
$ws2.Range("L96").Value = 4.985332489013672
$ws2.Range("M96").Value = 4.049539566040039
$ws2.Range("N96").Value = 9.034872055053711
$ws2.Range("M97").Value = 55.21059036254883
$ws2.Range("N97").Value = 55.21059036254883
$ws2.Range("M98").Value = 55.21059036254883
$ws2.Range("N98").Value = 55.21059036254883
$ws2.Range("L99").Value = 4.991292953491211
$ws2.Range("M99").Value = 5.012035369873047
$ws2.Range("N99").Value = 10.00332832336426
$ws2.Range("L100").Value = 4.991292953491211
$ws2.Range("M100").Value = 5.012035369873047
$ws2.Range("N100").Value = 10.00332832336426
$ws2.Range("L101").Value = 4.991292953491211
$ws2.Range("M101").Value = 5.012035369873047
$ws2.Range("N101").Value = 10.00332832336426
$ws2.Range("M102").Value = 10.00785827636719
$ws2.Range("N102").Value = 10.00785827636719
$ws2.Range("L103").Value = 12.98832893371582
$ws2.Range("N103").Value = 12.98832893371582
$ws2.Range("M104").Value = 5.960226058959961
$ws2.Range("N104").Value = 5.960226058959961
$ws2.Range("L105").Value = 15.98763465881348
$ws2.Range("M105").Value = 1.000404357910156
$ws2.Range("N105").Value = 16.98803901672363
$ws2.Range("M106").Value = 39.22796249389648
$ws2.Range("N106").Value = 39.22796249389648
$ws2.Range("L107").Value = 4.997491836547852
$ws2.Range("M107").Value = 2.999305725097656
$ws2.Range("N107").Value = 7.996797561645508
$ws2.Range("M108").Value = 35.37106513977051
$ws2.Range("N108").Value = 35.37106513977051
$ws2.Range("L109").Value = 5.545377731323242
$ws2.Range("M109").Value = 2.999067306518555
$ws2.Range("N109").Value = 8.544445037841797
$ws2.Range("L110").Value = 5.545377731323242
$ws2.Range("M110").Value = 2.999067306518555
$ws2.Range("N110").Value = 8.544445037841797
$ws2.Range("M111").Value = 21.99888229370117
$ws2.Range("N111").Value = 21.99888229370117
$ws2.Range("L112").Value = 7.001638412475586
$ws2.Range("N112").Value = 7.001638412475586
$ws2.Range("M113").Value = 34.03639793395996
$ws2.Range("N113").Value = 34.03639793395996
$ws2.Range("L114").Value = 20.01690864562988
$ws2.Range("N114").Value = 20.01690864562988
$ws2.Range("L115").Value = 0
$ws2.Range("M115").Value = 28.9008617401123
$ws2.Range("N115").Value = 28.9008617401123
$ws2.Range("L116").Value = 22.11570739746094
$ws2.Range("M116").Value = 3.004074096679688
$ws2.Range("N116").Value = 25.11978149414062
$ws2.Range("L117").Value = 3.986120223999023
$ws2.Range("M117").Value = 3.000020980834961
$ws2.Range("N117").Value = 6.986141204833984
$ws2.Range("L118").Value = 3.986120223999023
$ws2.Range("M118").Value = 3.000020980834961
$ws2.Range("N118").Value = 6.986141204833984
$ws2.Range("L119").Value = 6.667375564575195
$ws2.Range("M119").Value = 5.9967041015625
$ws2.Range("N119").Value = 12.6640796661377
$ws2.Range("L120").Value = 6.667375564575195
$ws2.Range("M120").Value = 5.9967041015625
$ws2.Range("N120").Value = 12.6640796661377
$ws2.Range("L121").Value = 9.058713912963867
$ws2.Range("M121").Value = 9.058475494384766
$ws2.Range("N121").Value = 18.11718940734863
$ws2.Range("L122").Value = 9.058713912963867
$ws2.Range("M122").Value = 9.058475494384766
$ws2.Range("N122").Value = 18.11718940734863
$ws2.Range("L123").Value = 17.22240447998047
$ws2.Range("M123").Value = 27.70376205444336
$ws2.Range("N123").Value = 44.92616653442383
$ws2.Range("L124").Value = 17.22240447998047
$ws2.Range("M124").Value = 27.70376205444336
$ws2.Range("N124").Value = 44.92616653442383
$ws2.Range("M125").Value = 187.4191761016846
$ws2.Range("N125").Value = 187.4191761016846
$ws2.Range("M126").Value = 187.4191761016846
$ws2.Range("N126").Value = 187.4191761016846
$ws2.Range("M127").Value = 187.4191761016846
$ws2.Range("N127").Value = 187.4191761016846
$ws2.Range("M128").Value = 187.4191761016846
$ws2.Range("N128").Value = 187.4191761016846
$ws2.Range("M129").Value = 187.4191761016846
$ws2.Range("N129").Value = 187.4191761016846
$ws2.Range("M130").Value = 187.4191761016846
$ws2.Range("N130").Value = 187.4191761016846
$ws2.Range("L131").Value = 24.48391914367676
$ws2.Range("M131").Value = 67.46840476989746
$ws2.Range("N131").Value = 91.95232391357422
$ws2.Range("L132").Value = 24.48391914367676
$ws2.Range("M132").Value = 67.46840476989746
$ws2.Range("N132").Value = 91.95232391357422
$ws2.Range("L133").Value = 24.48391914367676
$ws2.Range("M133").Value = 67.46840476989746
$ws2.Range("N133").Value = 91.95232391357422
$ws2.Range("L134").Value = 24.48391914367676
$ws2.Range("M134").Value = 67.46840476989746
$ws2.Range("N134").Value = 91.95232391357422

$ws3 = $wb.Worksheets.Item("summary_by_diagnoser")
$ws3.Range("E2").Value = 37.78591632843018
$ws3.Range("E3").Value = 18.94248962402344

Write-Host "done"